# Applies the tracked changes described by the commit diff:
#  1. Merge the "  " + <bookmark _GoBack/> + "TPBL" runs near "TPGD / TPBL"
#     into a single run "  TPBL" (the _GoBack bookmark is removed from here).
#  2. Colour the first stand-alone "CIF" label (customer block) red (FF0000).
#  3. Colour the second stand-alone "CIF" label (guarantor block) red
#     (FF0000) and re-anchor the _GoBack bookmark around that run.
#  4. Update the footer's cached PAGE field result from "1" to "7".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "  " + bookmark(_GoBack) + "TPBL"  ->  single run "  TPBL"
# ---------------------------------------------------------------------
$tpblFind = $d.Content.Duplicate
$tpblFind.Find.Execute("  TPBL") | Out-Null
if ($tpblFind.Find.Found) {
    $tpblStart = $tpblFind.Start
    $tpblEnd = $tpblFind.End

    # Force a real text mutation (assigning the identical string is a no-op
    # in this host) so the bookmark straddling the run boundary is dropped,
    # then write the final merged text back into the same span.
    $tpblRange = $d.Range($tpblStart, $tpblEnd)
    $tpblRange.Text = [string][char]1
    $tpblRange2 = $d.Range($tpblStart, $tpblStart + 1)
    $tpblRange2.Text = "  TPBL"
}

# ---------------------------------------------------------------------
# 2) & 3) Colour the two "CIF" labels red; move _GoBack onto the 2nd one
# ---------------------------------------------------------------------
$cifHits = New-Object System.Collections.ArrayList
$cursor = $d.Content.Duplicate
$cursor.Start = 0
while ($cursor.Find.Execute("CIF")) {
    $cifHits.Add(@($cursor.Start, $cursor.End)) | Out-Null
    $cursor.Collapse(0)
    $cursor.End = $d.Content.End
}

for ($i = 0; $i -lt $cifHits.Count; $i++) {
    $pair = $cifHits[$i]
    $cifRange = $d.Range($pair[0], $pair[1])
    $cifRange.Font.Color = 255  # wdColorRed -> <w:color w:val="FF0000"/>

    if ($i -eq ($cifHits.Count - 1)) {
        # Re-create the _GoBack bookmark around the last "CIF" run.
        $d.Bookmarks.Add("_GoBack", $cifRange)
    }
}

# ---------------------------------------------------------------------
# 4) Footer cached PAGE field result: "1" -> "7"
# ---------------------------------------------------------------------
$footerStory = $d.StoryRanges(9)  # wdPrimaryFooterStory
$pageNumRange = $footerStory.Duplicate
$pageNumRange.Find.ClearFormatting()
$pageNumRange.Find.Execute("1") | Out-Null
if ($pageNumRange.Find.Found) {
    $pnStart = $pageNumRange.Start
    $pnEnd = $pageNumRange.End
    $target = $footerStory.Duplicate
    $target.Start = $pnStart
    $target.End = $pnEnd
    $target.Text = "7"
}
